$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting so that
# numeric-looking strings (e.g. "241.50", "0.6085") are not silently coerced
# into floating point numbers by Excels type inference.
$ws.Range("D2:E51").NumberFormat = "@"

$data = @(
    @{Row=2; B='Bitcoin'; C='https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D='29.011.91'; E='  -0.45%  '},
    @{Row=3; B='Ethereum'; C='https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D='1.819.54'; E='  -1.04%  '},
    @{Row=4; B='TetherUSD'; C='https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D='1.004'; E='  +0.20%  '},
    @{Row=5; B='BNB'; C='https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D='241.50'; E='  -1.22%  '},
    @{Row=6; B='XRP'; C='https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D='0.6085'; E='  -3.60%  '},
    @{Row=7; B='USDC'; C='https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D='1.005'; E='  +0.11%  '},
    @{Row=8; B='Dogecoin'; C='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D='0.07299'; E='  -2.91%  '},
    @{Row=9; B='Cardano'; C='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D='0.2867'; E='  -2.28%  '},
    @{Row=10; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='22.64'; E='  -2.69%  '},
    @{Row=11; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.07656'; E='  -1.14%  '},
    @{Row=12; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='1.831.97'; E='  -0.31%  '},
    @{Row=13; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='4.911'; E='  -1.77%  '},
    @{Row=14; B='Polygon'; C='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D='0.6544'; E='  -2.44%  '},
    @{Row=15; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='80.84'; E='  -2.33%  '},
    @{Row=16; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.000008858'; E='  -4.95%  '},
    @{Row=17; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='5.822'; E='  -3.37%  '},
    @{Row=18; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='29.014.89'; E='  -0.53%  '},
    @{Row=19; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='2.069.32'; E='  -0.59%  '},
    @{Row=20; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='234.40'; E='  +4.59%  '},
    @{Row=21; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='12.35'; E='  -2.08%  '},
    @{Row=22; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='1.005'; E='  +0.01%  '},
    @{Row=23; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='7.062'; E='  -1.28%  '},
    @{Row=24; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='1.005'; E='  +0.21%  '},
    @{Row=25; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='158.73'; E='  -0.80%  '},
    @{Row=26; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.1385'; E='  -1.33%  '},
    @{Row=27; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='8.369'; E='  -1.87%  '},
    @{Row=28; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='17.53'; E='  -2.55%  '},
    @{Row=29; B='PancakeSwap'; C='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D='1.486'; E='  -1.25%  '},
    @{Row=30; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.05577'; E='  -5.66%  '},
    @{Row=31; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='4.046'; E='  -0.69%  '},
    @{Row=32; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='1.207'; E='  -0.02%  '},
    @{Row=33; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='4.048'; E='  -2.85%  '},
    @{Row=34; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='1.814'; E='  -2.15%  '},
    @{Row=35; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.7235'; E='  -3.67%  '},
    @{Row=36; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='1.126'; E='  -1.57%  '},
    @{Row=37; B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='2.632'; E='  -1.80%  '},
    @{Row=38; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='2.806'; E='  +1.17%  '},
    @{Row=39; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.01750'; E='  -2.45%  '},
    @{Row=40; B='Maker'; C='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D='1.190.79'; E='  -3.50%  '},
    @{Row=41; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='6.335'; E='  -3.72%  '},
    @{Row=42; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='0.8811'; E='  -1.70%  '},
    @{Row=43; B='PaxDollar'; C='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D='1.005'; E='  +0.03%  '},
    @{Row=44; B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='100.67'; E='  -1.73%  '},
    @{Row=45; B='RocketPoolETH'; C='https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'; D='1.979.29'; E='  +0.00%  '},
    @{Row=46; B='BabyDogeCoin'; C='https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; D='0.00000000121'; E='  -2.67%  '},
    @{Row=47; B='Aave'; C='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D='63.91'; E='  -3.25%  '},
    @{Row=48; B='Mantle'; C='https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; D='0.5101'; E='  +0.01%  '},
    @{Row=49; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='9.007'; E='  -0.37%  '},
    @{Row=50; B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.3959'; E='  -2.98%  '},
    @{Row=51; B='Cronos'; C='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D='0.05787'; E='  -0.85%  '}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
    $ws.Cells.Item($item.Row, 4).Value = $item.D
    $ws.Cells.Item($item.Row, 5).Value = $item.E
}
